# Apply is_word check to lexical density and lexical diversity
# This updates the per-document lexical_density values (sheets 2019-2024)
# and the derived Summary statistics (sheet "Summary") to reflect the
# recomputed values after filtering tokens with an is_word() check.

$wb = $excel.ActiveWorkbook

# --- Sheet "2019": column B lexical_density values (rows 2..5) ---
$ws = $wb.Worksheets.Item("2019")
$ws.Range("B2").Value = 0.5833333333333334
$ws.Range("B3").Value = 0.5107398568019093
$ws.Range("B4").Value = 0.555045871559633
$ws.Range("B5").Value = 0.5818033455732354

# --- Sheet "2020": column B lexical_density values (rows 2..8) ---
$ws = $wb.Worksheets.Item("2020")
$ws.Range("B2").Value = 0.6496350364963503
$ws.Range("B3").Value = 0.5366568914956011
$ws.Range("B4").Value = 0.6271929824561403
$ws.Range("B5").Value = 0.625
$ws.Range("B6").Value = 0.5308641975308642
$ws.Range("B7").Value = 0.5644067796610169
$ws.Range("B8").Value = 0.6055555555555555

# --- Sheet "2021": column B lexical_density values (rows 2..11) ---
$ws = $wb.Worksheets.Item("2021")
$ws.Range("B2").Value = 0.5436893203883495
$ws.Range("B3").Value = 0.589041095890411
$ws.Range("B4").Value = 0.5705128205128205
$ws.Range("B5").Value = 0.559322033898305
$ws.Range("B6").Value = 0.6320166320166321
$ws.Range("B7").Value = 0.5454545454545454
$ws.Range("B8").Value = 0.5390070921985816
$ws.Range("B9").Value = 0.6151832460732984
$ws.Range("B10").Value = 0.5789473684210527
$ws.Range("B11").Value = 0.5857142857142857

# --- Sheet "2022": column B lexical_density values (rows 2..7) ---
$ws = $wb.Worksheets.Item("2022")
$ws.Range("B2").Value = 0.5955882352941176
$ws.Range("B3").Value = 0.5277337559429477
$ws.Range("B4").Value = 0.5365853658536586
$ws.Range("B5").Value = 0.6024464831804281
$ws.Range("B6").Value = 0.5643115942028986
$ws.Range("B7").Value = 0.5587734241908007

# --- Sheet "2023": column B lexical_density values (row 2) ---
$ws = $wb.Worksheets.Item("2023")
$ws.Range("B2").Value = 0.565597667638484

# --- Sheet "2024": column B lexical_density values (rows 2..4) ---
$ws = $wb.Worksheets.Item("2024")
$ws.Range("B2").Value = 0.6065573770491803
$ws.Range("B3").Value = 0.5759717314487632
$ws.Range("B4").Value = 0.5602409638554217

# --- Sheet "Summary": columns C..I (mean, std, min, 25%, 50%, 75%, max) ---
$ws = $wb.Worksheets.Item("Summary")

# Row 2 -> year 2019
$ws.Range("C2").Value = 0.5577306018170278
$ws.Range("D2").Value = 0.03391329506589277
$ws.Range("E2").Value = 0.5107398568019093
$ws.Range("F2").Value = 0.5439693678702021
$ws.Range("G2").Value = 0.5684246085664342
$ws.Range("H2").Value = 0.5821858425132599
$ws.Range("I2").Value = 0.5833333333333334

# Row 3 -> year 2020
$ws.Range("C3").Value = 0.5913302061707898
$ws.Range("D3").Value = 0.04724367325764119
$ws.Range("E3").Value = 0.5308641975308642
$ws.Range("F3").Value = 0.5505318355783091
$ws.Range("G3").Value = 0.6055555555555555
$ws.Range("H3").Value = 0.6260964912280702
$ws.Range("I3").Value = 0.6496350364963503

# Row 4 -> year 2021
$ws.Range("C4").Value = 0.5758888440568282
$ws.Range("D4").Value = 0.03091194396711534
$ws.Range("E4").Value = 0.5390070921985816
$ws.Range("F4").Value = 0.5489214175654853
$ws.Range("G4").Value = 0.5747300944669366
$ws.Range("H4").Value = 0.5882093933463797
$ws.Range("I4").Value = 0.6320166320166321

# Row 5 -> year 2022
$ws.Range("C5").Value = 0.5642398097774752
$ws.Range("D5").Value = 0.03023183896021963
$ws.Range("E5").Value = 0.5277337559429477
$ws.Range("F5").Value = 0.5421323804379441
$ws.Range("G5").Value = 0.5615425091968496
$ws.Range("H5").Value = 0.5877690750213129
$ws.Range("I5").Value = 0.6024464831804281

# Row 6 -> year 2023 (single sample; no std/D6)
$ws.Range("C6").Value = 0.565597667638484
$ws.Range("E6").Value = 0.565597667638484
$ws.Range("F6").Value = 0.565597667638484
$ws.Range("G6").Value = 0.565597667638484
$ws.Range("H6").Value = 0.565597667638484
$ws.Range("I6").Value = 0.565597667638484

# Row 7 -> year 2024
$ws.Range("C7").Value = 0.5809233574511218
$ws.Range("D7").Value = 0.02355188915637492
$ws.Range("E7").Value = 0.5602409638554217
$ws.Range("F7").Value = 0.5681063476520924
$ws.Range("G7").Value = 0.5759717314487632
$ws.Range("H7").Value = 0.5912645542489718
$ws.Range("I7").Value = 0.6065573770491803
